# [MINDX] MODULE 2 - BUOI 11
#
# Insert a new "Hr Attrtion Dashboard" divider slide as the 2nd slide of
# the deck (it becomes sldId 259, sitting between the existing sldId 256
# and sldId 258 slides). The new slide is built by duplicating slide 1
# (which already carries the shared header/icon/title treatment used by
# this deck) and then adjusting its shape set so it matches the new
# slide's content: drop the two small stat cards and add a single large
# rounded-rectangle card below the title card.

$p = $ppt.ActivePresentation

# --- 1. Duplicate slide 1 into slide position 2 -----------------------
$s1 = $p.Slides.Item(1)
$s1.Duplicate() | Out-Null
$s2 = $p.Slides.Item(2)

# --- 2. Remove the two small stat cards that don't belong on the new --
#        slide ("Rounded Rectangle 6" / "Rounded Rectangle 7").
$s2.Shapes.Item("Rounded Rectangle 6").Delete()
$s2.Shapes.Item("Rounded Rectangle 7").Delete()

# --- 3. Add the single big card shape that replaces them --------------
# Shape position/size on the COM surface is expressed in points
# (1 pt = 12700 EMU); the target geometry in EMU is
#   off  x=4692650  y=3556000
#   ext cx=6995160 cy=2862072
$cardLeft   = 4692650 / 12700.0
$cardTop    = 3556000 / 12700.0
$cardWidth  = 6995160 / 12700.0
$cardHeight = 2862072 / 12700.0

$card = $s2.Shapes.AddShape(5, $cardLeft, $cardTop, $cardWidth, $cardHeight)
$card.Name = "Rounded Rectangle 1"
$card.Adjustments.Item(1) = 0.09784
$card.Fill.ForeColor.RGB = 0x3C2F2E
$card.Line.Visible = $false
$card.TextFrame.VerticalAnchor = 3
$card.TextFrame.TextRange.ParagraphFormat.Alignment = 2
